# Update the "取得日時" (retrieved-at) timestamp in column A for the
# existing data rows (rows 2-11) on the active sheet ("ランサーズ")
# to reflect the latest scrape run: 2025-10-21 12:49:33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-10-21 12:49:33"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 11
}

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
